$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.883.03"
$ws.Range("E2").Value = "  -0.91%  "
# Row 3
$ws.Range("D3").Value = "2.367.93"
$ws.Range("E3").Value = "  -1.44%  "
# Row 4
$ws.Range("E4").Value = "  +0.05%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.90%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.38%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.59%  "
# Row 8
$ws.Range("E8").Value = "  +0.03%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.51%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.22%  "
# Row 11
$ws.Range("E11").Value = "  -1.47%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.96%  "
# Row 13
$ws.Range("E13").Value = "  -4.03%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.76%  "
# Row 16
$ws.Range("D16").Value = "2.726.30"
$ws.Range("E16").Value = "  -1.38%  "
# Row 17
$ws.Range("D17").Value = "2.348.51"
$ws.Range("E17").Value = "  -1.94%  "
# Row 18
$ws.Range("D18").Value = "42.848.36"
$ws.Range("E18").Value = "  -0.97%  "
# Row 19
$ws.Range("E19").Value = "  -0.58%  "
# Row 20
$ws.Range("E20").Value = "  -2.03%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.14%  "
# Row 22
$ws.Range("E22").Value = "  -0.62%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.94%  "
# Row 24
$ws.Range("E24").Value = "  -3.82%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.96%  "
# Row 26
$ws.Range("E26").Value = "  +0.06%  "
# Row 27
$ws.Range("E27").Value = "  -2.83%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.89%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.41"
$ws.Range("D30").Style = "Normal"
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.20%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0898"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.73%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.95%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.89%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.121"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.20%  "
# Row 36
$ws.Range("E36").Value = "  -2.44%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.97%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0367"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.21%  "
# Row 40
$ws.Range("E40").Value = "  -5.23%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.243"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.61%  "
# Row 42
$ws.Range("E42").Value = "  -4.64%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.05%  "
# Row 44
$ws.Range("E44").Value = "  +0.02%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "
# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "113.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.66%  "
# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.31%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.98%  "
# Row 51
$ws.Range("E51").Value = "  -0.84%  "
